$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing data row (374) down through the new rows (375-385)
$ws.Range("A374:D374").Copy()
$ws.Range("A375:D385").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new rows with data (aggiornamento fino a 20/09/2021)
$ws.Cells.Item(375, 1).Value = 44449
$ws.Cells.Item(375, 2).Value = 1
$ws.Cells.Item(375, 3).Value = 27
$ws.Cells.Item(375, 4).Value = 150.1835576816109
$ws.Cells.Item(376, 1).Value = 44450
$ws.Cells.Item(376, 2).Value = 5
$ws.Cells.Item(376, 3).Value = 30
$ws.Cells.Item(376, 4).Value = 166.8706196462343
$ws.Cells.Item(377, 1).Value = 44451
$ws.Cells.Item(377, 2).Value = 5
$ws.Cells.Item(377, 3).Value = 32
$ws.Cells.Item(377, 4).Value = 177.9953276226499
$ws.Cells.Item(378, 1).Value = 44452
$ws.Cells.Item(378, 2).Value = 1
$ws.Cells.Item(378, 3).Value = 28
$ws.Cells.Item(378, 4).Value = 155.7459116698187
$ws.Cells.Item(379, 1).Value = 44453
$ws.Cells.Item(379, 2).Value = 4
$ws.Cells.Item(379, 3).Value = 31
$ws.Cells.Item(379, 4).Value = 172.4329736344421
$ws.Cells.Item(380, 1).Value = 44454
$ws.Cells.Item(380, 2).Value = 0
$ws.Cells.Item(380, 3).Value = 31
$ws.Cells.Item(380, 4).Value = 172.4329736344421
$ws.Cells.Item(381, 1).Value = 44455
$ws.Cells.Item(381, 2).Value = 2
$ws.Cells.Item(381, 3).Value = 18
$ws.Cells.Item(381, 4).Value = 100.1223717877406
$ws.Cells.Item(382, 1).Value = 44456
$ws.Cells.Item(382, 2).Value = 8
$ws.Cells.Item(382, 3).Value = 25
$ws.Cells.Item(382, 4).Value = 139.0588497051953
$ws.Cells.Item(383, 1).Value = 44457
$ws.Cells.Item(383, 2).Value = 2
$ws.Cells.Item(383, 3).Value = 22
$ws.Cells.Item(383, 4).Value = 122.3717877405718
$ws.Cells.Item(384, 1).Value = 44458
$ws.Cells.Item(384, 2).Value = 4
$ws.Cells.Item(384, 3).Value = 21
$ws.Cells.Item(384, 4).Value = 116.809433752364
$ws.Cells.Item(385, 1).Value = 44459
$ws.Cells.Item(385, 2).Value = 0
$ws.Cells.Item(385, 3).Value = 20
$ws.Cells.Item(385, 4).Value = 111.2470797641562
